# TC24_Canine_Filter_Breed-FrenchBullDg.xlsx
# "Fixed variables and query errors in Bread from TC01 to TC30"
#
# The CasesTab Cypher query (startup!B2) referenced a `co.cohort_description`
# / `Cohort` column that was dropped from the final RETURN clause. Update the
# cell in place so the shared-string table, the row/column references, and
# the rest of the workbook stay internally consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")
$ws.Activate()

$fixedQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n" +
  "WHERE demo.breed IN ['French Bulldog']`n" +
  "MATCH (c)<--(diag:diagnosis)`n" +
  "OPTIONAL MATCH (samp:sample)-->(c)`n" +
  "OPTIONAL MATCH (co:cohort)<-[*]-(c)`n" +
  "WITH DISTINCT c, s, demo, diag, co`n" +
  "RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n" +
  "        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n" +
  "        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n" +
  "        coalesce(demo.breed, '') AS Breed ,`n" +
  "        coalesce(diag.disease_term, '') AS Diagnosis ,`n" +
  "        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n" +
  "        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n" +
  "        coalesce(demo.sex, '') AS Sex ,`n" +
  "        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n" +
  "        coalesce(demo.weight, '') AS ``Weight (kg)``,`n" +
  "        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $fixedQuery

# The shorter query text wraps to fewer lines; match the new autofit height.
$ws.Rows.Item(2).RowHeight = 244.8

# Leave the edited cell selected (matches the author's last-saved selection).
$ws.Range("B2").Select() | Out-Null
